$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new outcome measurements for the "Pre Experimental Phase" (column C)
$ws.Range("C2").Value = "A little stressful"
$ws.Range("C3").Value = "A little stressful"
$ws.Range("C4").Value = "Not stressful"
$ws.Range("C5").Value = "Not stressful"
$ws.Range("C6").Value = "Moderately stressful"
$ws.Range("C7").Value = "Moderately stressful"

# Update selected cell to reflect new data entry position
$ws.Range("C8").Select()
